# Add "Estatus" column (F) to the BACKLOG table with per-row OK/NOOK values,
# resize the table, add conditional formatting, and tidy up the view/column
# width to match the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Resize the table "Tabla1" from B4:E26 to B4:F26 (adds 5th column) ---
$lo = $ws.ListObjects.Item("Tabla1")
$lo.Resize($ws.Range("B4:F26"))

# --- 2. Header + per-row Estatus values (this also names the new table column) ---
$ws.Range("F4").Value = "Estatus"

$ws.Range("F5").Value  = "OK"
$ws.Range("F6").Value  = "NOOK"
$ws.Range("F7").Value  = "NOOK"
$ws.Range("F8").Value  = "NOOK"
$ws.Range("F9").Value  = "NOOK"
$ws.Range("F10").Value = "OK"
$ws.Range("F11").Value = "NOOK"
$ws.Range("F12").Value = "NOOK"
$ws.Range("F13").Value = "NOOK"
$ws.Range("F14").Value = "NOOK"
$ws.Range("F15").Value = "NOOK"
$ws.Range("F16").Value = "OK"
$ws.Range("F17").Value = "OK"
$ws.Range("F18").Value = "OK"
$ws.Range("F19").Value = "NOOK"
$ws.Range("F20").Value = "OK"

# --- 3. Conditional formatting: OK -> green fill, NOOK -> red fill (bold white text) ---
$rng = $ws.Range("F5:F26")

$fcOk = $rng.FormatConditions.Add(1, 3, '"OK"')
$fcOk.Font.Bold = $true
$fcOk.Font.Italic = $false
$fcOk.Interior.Color = 5296274
$fcOk.Priority = 3

$fcNook = $rng.FormatConditions.Add(1, 3, '"NOOK"')
$fcNook.Font.Bold = $true
$fcNook.Font.Italic = $false
$fcNook.Interior.Color = 255
$fcNook.Priority = 1

# --- 4. Widen column F to fit the new header text ---
$ws.Columns("F").ColumnWidth = 7.29

# --- 5. Merge the title row across the new column too ---
$ws.Range("B2:E2").UnMerge()
$ws.Range("B2:F2").Merge()
$ws.Range("B2:F2").Borders.LineStyle = 0
$ws.Range("B2").Borders(7).LineStyle = 1

# --- 6. Reset view: active cell back to C4, no saved scroll position ---
$ws.Range("C4").Select()
